$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 927.2727
$ws.Range("J45").Value = 900
$ws.Range("L45").Value = 900
$ws.Range("N45").Value = -1654
$ws.Range("H63").Value = 771228.0600000001
$ws.Range("I63").Value = 911087.75
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 911087.75
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -910401.75
$ws.Range("N63").Value = -3372
$ws.Range("H66").Value = 771228.0600000001
$ws.Range("I66").Value = 911087.75
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 4555438.75
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -4552006.75
$ws.Range("N66").Value = -16864
$ws.Range("H110").Value = 1167.2693
$ws.Range("I110").Value = 847.5789
$ws.Range("K110").Value = 847.5789
$ws.Range("M110").Value = 1197.4211
$ws.Range("H123").Value = 51476
$ws.Range("J123").Value = 51476
$ws.Range("L123").Value = 51476
$ws.Range("N123").Value = -61276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2044.0526
$ws.Range("I20").Value = 3229.875
$ws.Range("J20").Value = 1181.6364
$ws.Range("K20").Value = 3229.875
$ws.Range("L20").Value = 1181.6364
$ws.Range("M20").Value = -2982.875
$ws.Range("N20").Value = -1675.6364

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5558651.5
$ws.Range("I31").Value = 2165.2727
$ws.Range("K31").Value = 2165.2727
$ws.Range("M31").Value = -1870.2727
$ws.Range("H34").Value = 5558651.5
$ws.Range("I34").Value = 2165.2727
$ws.Range("K34").Value = 2165.2727
$ws.Range("M34").Value = -1963.2727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 942.8570999999999
$ws.Range("I5").Value = 561.25
$ws.Range("J5").Value = 1451.6666
$ws.Range("K5").Value = 1683.75
$ws.Range("L5").Value = 4354.9998
$ws.Range("M5").Value = -1571.75
$ws.Range("N5").Value = -4578.9998
$ws.Range("H63").Value = 7963.3335
$ws.Range("I63").Value = 5260
$ws.Range("J63").Value = 10666.667
$ws.Range("K63").Value = 15780
$ws.Range("L63").Value = 32000.001
$ws.Range("M63").Value = -15031
$ws.Range("N63").Value = -33498.001
$ws.Range("H64").Value = 1887.5883
$ws.Range("I64").Value = 972.25
$ws.Range("J64").Value = 2169.2307
$ws.Range("K64").Value = 2916.75
$ws.Range("L64").Value = 6507.6921
$ws.Range("M64").Value = -2646.75
$ws.Range("N64").Value = -7047.6921
$ws.Range("H66").Value = 7963.3335
$ws.Range("I66").Value = 5260
$ws.Range("J66").Value = 10666.667
$ws.Range("K66").Value = 47340
$ws.Range("L66").Value = 96000.003
$ws.Range("M66").Value = -43596
$ws.Range("N66").Value = -103488.003
$ws.Range("H67").Value = 1887.5883
$ws.Range("I67").Value = 972.25
$ws.Range("J67").Value = 2169.2307
$ws.Range("K67").Value = 2916.75
$ws.Range("L67").Value = 6507.6921
$ws.Range("M67").Value = -1980.75
$ws.Range("N67").Value = -8379.6921
$ws.Range("H87").Value = 979.6667
$ws.Range("I87").Value = 979.6667
$ws.Range("K87").Value = 2939.0001
$ws.Range("M87").Value = -1691.0001
$ws.Range("H90").Value = 979.6667
$ws.Range("I90").Value = 979.6667
$ws.Range("K90").Value = 8817.0003
$ws.Range("M90").Value = -2577.0003
$ws.Range("H114").Value = 4724.4614
$ws.Range("I114").Value = 1746.8572
$ws.Range("J114").Value = 8198.333000000001
$ws.Range("K114").Value = 5240.571599999999
$ws.Range("L114").Value = 24594.999
$ws.Range("M114").Value = -1986.571599999999
$ws.Range("N114").Value = -31102.999
$ws.Range("H122").Value = 1468.6333
$ws.Range("I122").Value = 1590.1818
$ws.Range("K122").Value = 14311.6362
$ws.Range("M122").Value = -11861.6362
$ws.Range("H132").Value = 1512.4375
$ws.Range("I132").Value = 528.4286
$ws.Range("J132").Value = 2277.7778
$ws.Range("K132").Value = 4755.8574
$ws.Range("L132").Value = 20500.0002
$ws.Range("M132").Value = -2225.8574
$ws.Range("N132").Value = -25560.0002
$ws.Range("H135").Value = 942.8570999999999
$ws.Range("I135").Value = 561.25
$ws.Range("J135").Value = 1451.6666
$ws.Range("K135").Value = 5051.25
$ws.Range("L135").Value = 13064.9994
$ws.Range("M135").Value = -2516.25
$ws.Range("N135").Value = -18134.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8504789
$ws.Range("I70").Value = 10855502
$ws.Range("K70").Value = 10855502
$ws.Range("M70").Value = -10855232
$ws.Range("H73").Value = 8504789
$ws.Range("I73").Value = 10855502
$ws.Range("K73").Value = 10855502
$ws.Range("M73").Value = -10854566
$ws.Range("H135").Value = 63450
$ws.Range("J135").Value = 63450
$ws.Range("L135").Value = 63450
$ws.Range("N135").Value = -73590

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2709.0908
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 2960
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 2960
$ws.Range("M40").Value = -2364
$ws.Range("N40").Value = -3232
$ws.Range("H61").Value = 1994.909
$ws.Range("I61").Value = 1183.1666
$ws.Range("J61").Value = 2969
$ws.Range("K61").Value = 1183.1666
$ws.Range("L61").Value = 2969
$ws.Range("M61").Value = -981.1666
$ws.Range("N61").Value = -3373
$ws.Range("H113").Value = 1994.909
$ws.Range("I113").Value = 1183.1666
$ws.Range("J113").Value = 2969
$ws.Range("K113").Value = 1183.1666
$ws.Range("L113").Value = 2969
$ws.Range("M113").Value = 986.8334
$ws.Range("N113").Value = -7309

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 314.81818
$ws.Range("I107").Value = 273.66666
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 820.9999799999999
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 1099.00002
$ws.Range("N107").Value = -5340
$ws.Range("H113").Value = 615.5
$ws.Range("I113").Value = 801.7
$ws.Range("J113").Value = 150
$ws.Range("K113").Value = 2405.1
$ws.Range("L113").Value = 450
$ws.Range("M113").Value = -235.1000000000004
$ws.Range("N113").Value = -4790
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
